# daily auto push: 2026-01-26 09:47 UTC
#
# The source data table (A:D = 日付/曜日/時刻/ランキング) is missing one
# observation for 2026/01/26 (月). Insert a new row right after the last
# existing 2026/01/26 row (row 717) carrying the next 時刻 value (16),
# which pushes every following row down by one and grows the used range
# from A1:D759 to A1:D760.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the row immediately above the insertion point (same date /
# weekday text, already stored as plain text - not as an auto-parsed
# date) and insert the copy as a new row 718. This shifts rows 718..759
# down to 719..760 and keeps cell formatting/types identical to their
# neighbours (avoids Excel's auto-detection turning "2026/01/26" into a
# real date value when typed in fresh).
$ws.Range("A717:D717").Copy()
$ws.Rows.Item(718).Insert()

# Only the 時刻 (hour) column differs for the newly added observation.
$ws.Range("C718").Value = 16
